$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.888.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.021.80'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.05%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.24'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.24'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +8.44%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.011.10'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.135'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.32'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +11.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.461'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000233'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.26'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.26%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.520.79'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.22'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.020.50'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '59.860.00'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '439.81'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +5.52%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.723'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.15'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.37'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.87'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.25'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +13.36%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.87'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.10'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.95%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.31%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +7.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0797'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +17.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.94'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.19%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.23'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.63%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.82'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +10.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '406.56'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.68%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0355'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.785.11'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.74%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.255'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.00%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.88'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.69%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.99'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +20.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.70'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.93%  '
